$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 15
$ws.Range("M15").Value = 1
$ws.Range("N15").Value = 2

# Row 31
$ws.Range("M31").Value = 3
$ws.Range("N31").Value = 3

# Row 32
$ws.Range("M32").Value = 1

# Row 33
$ws.Range("M33").Value = 3
$ws.Range("N33").Value = 3

# Row 39
$ws.Range("M39").Value = 2
$ws.Range("N39").Value = 2

# Row 42
$ws.Range("M42").Value = 4

# Row 49
$ws.Range("M49").Value = 2
$ws.Range("N49").Value = 3
